$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# "Tags" block (rows 12-14): harmonize the term-accession / term-source tags
# for the "Mass Spectrometry" / "data processing protocol" columns so they
# use the same short CURIE style as the rest of the sheet.
$ws.Range("C13").Value = "NCIT:C17156"
$ws.Range("D12").Value = "data processing"
$ws.Range("D13").Value = "NCIT:C47925"
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()

# Row 13 used a manually-wrapped, auto-fit height to fit the long URL that
# used to live in C13/D13; now that the text is short it no longer needs to
# wrap, so recompute the row height.
$ws.Rows.Item(13).EntireRow.AutoFit()

# Move the active selection to B17 on the isa_template sheet.
$ws.Activate()
$ws.Range("B17").Select()
